# Auto-generated: apply cell-value updates per the target diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 3075.111
$ws.Range("I111").Value = 2255.4285
$ws.Range("K111").Value = 6766.2855
$ws.Range("M111").Value = -3699.2855
$ws.Range("H125").Value = 3472
$ws.Range("I125").Value = 1081
$ws.Range("K125").Value = 9729
$ws.Range("M125").Value = -7269
$ws.Range("H132").Value = 14119.5
$ws.Range("I132").Value = 2283.4595
$ws.Range("K132").Value = 6850.3785
$ws.Range("M132").Value = -4320.3785
$ws.Range("H135").Value = 29414488
$ws.Range("I135").Value = 31252768
$ws.Range("K135").Value = 281274912
$ws.Range("M135").Value = -281272377
$ws.Range("H137").Value = 2469.4695
$ws.Range("I137").Value = 2333.6667
$ws.Range("J137").Value = 3284.2856
$ws.Range("K137").Value = 7001.000100000001
$ws.Range("L137").Value = 9852.856800000001
$ws.Range("M137").Value = -4451.000100000001
$ws.Range("N137").Value = -14952.8568
$ws.Range("H141").Value = 5024.567
$ws.Range("I141").Value = 2951.8462
$ws.Range("K141").Value = 8855.5386
$ws.Range("M141").Value = -3675.5386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1357.73
$ws.Range("I32").Value = 1367.7858
$ws.Range("K32").Value = 1367.7858
$ws.Range("M32").Value = -1080.7858
$ws.Range("H38").Value = 608
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H45").Value = 3354.6667
$ws.Range("I45").Value = 2350
$ws.Range("K45").Value = 2350
$ws.Range("M45").Value = -1973
$ws.Range("H61").Value = 2602.2903
$ws.Range("I61").Value = 2454.28
$ws.Range("K61").Value = 2454.28
$ws.Range("M61").Value = -2242.28
$ws.Range("H74").Value = 2772.6316
$ws.Range("I74").Value = 3210.6924
$ws.Range("J74").Value = 1823.5
$ws.Range("K74").Value = 3210.6924
$ws.Range("L74").Value = 1823.5
$ws.Range("M74").Value = -2336.6924
$ws.Range("N74").Value = -3571.5
$ws.Range("H77").Value = 2772.6316
$ws.Range("I77").Value = 3210.6924
$ws.Range("J77").Value = 1823.5
$ws.Range("K77").Value = 16053.462
$ws.Range("L77").Value = 9117.5
$ws.Range("M77").Value = -11685.462
$ws.Range("N77").Value = -17853.5
$ws.Range("H97").Value = 1292.0513
$ws.Range("I97").Value = 775.1515000000001
$ws.Range("J97").Value = 4135
$ws.Range("K97").Value = 775.1515000000001
$ws.Range("L97").Value = 4135
$ws.Range("M97").Value = -279.1515000000001
$ws.Range("N97").Value = -5127
$ws.Range("H122").Value = 2161.4375
$ws.Range("I122").Value = 1608.6
$ws.Range("K122").Value = 4825.799999999999
$ws.Range("M122").Value = -2375.799999999999
$ws.Range("H132").Value = 2138.0938
$ws.Range("I132").Value = 1735.2609
$ws.Range("J132").Value = 3167.5557
$ws.Range("K132").Value = 5205.7827
$ws.Range("L132").Value = 9502.667099999999
$ws.Range("M132").Value = -2675.7827
$ws.Range("N132").Value = -14562.6671
$ws.Range("H136").Value = 2602.2903
$ws.Range("I136").Value = 2454.28
$ws.Range("K136").Value = 7362.84
$ws.Range("M136").Value = -4812.84

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4173.303
$ws.Range("I86").Value = 3416.88
$ws.Range("J86").Value = 6537.125
$ws.Range("K86").Value = 3416.88
$ws.Range("L86").Value = 6537.125
$ws.Range("M86").Value = -2293.88
$ws.Range("N86").Value = -8783.125
$ws.Range("H89").Value = 4173.303
$ws.Range("I89").Value = 3416.88
$ws.Range("J89").Value = 6537.125
$ws.Range("K89").Value = 17084.4
$ws.Range("L89").Value = 32685.625
$ws.Range("M89").Value = -11468.4
$ws.Range("N89").Value = -43917.625
$ws.Range("H94").Value = 2056.2104
$ws.Range("I94").Value = 2049.0344
$ws.Range("J94").Value = 2079.3333
$ws.Range("K94").Value = 2049.0344
$ws.Range("L94").Value = 2079.3333
$ws.Range("M94").Value = -1598.0344
$ws.Range("N94").Value = -2981.3333
$ws.Range("H99").Value = 32846.383
$ws.Range("I99").Value = 46367
$ws.Range("K99").Value = 46367
$ws.Range("M99").Value = -44869
$ws.Range("H107").Value = 12467.75
$ws.Range("I107").Value = 9575.691999999999
$ws.Range("K107").Value = 9575.691999999999
$ws.Range("M107").Value = -7655.691999999999
$ws.Range("H133").Value = 60000
$ws.Range("J133").Value = 60000
$ws.Range("L133").Value = 60000
$ws.Range("N133").Value = -70120
$ws.Range("H134").Value = 1770.8334
$ws.Range("I134").Value = 1279.8914
$ws.Range("K134").Value = 3839.6742
$ws.Range("M134").Value = -1304.6742

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2180.8235
$ws.Range("I31").Value = 2069.2964
$ws.Range("J31").Value = 2611
$ws.Range("K31").Value = 2069.2964
$ws.Range("L31").Value = 2611
$ws.Range("M31").Value = -1774.2964
$ws.Range("N31").Value = -3201
$ws.Range("H34").Value = 2180.8235
$ws.Range("I34").Value = 2069.2964
$ws.Range("J34").Value = 2611
$ws.Range("K34").Value = 2069.2964
$ws.Range("L34").Value = 2611
$ws.Range("M34").Value = -1867.2964
$ws.Range("N34").Value = -3015
$ws.Range("H58").Value = 2349.32
$ws.Range("I58").Value = 1402.1666
$ws.Range("K58").Value = 1402.1666
$ws.Range("M58").Value = -1199.1666
$ws.Range("H105").Value = 4731.357
$ws.Range("I105").Value = 4228.25
$ws.Range("J105").Value = 7750
$ws.Range("K105").Value = 4228.25
$ws.Range("L105").Value = 7750
$ws.Range("M105").Value = -2481.25
$ws.Range("N105").Value = -11244
$ws.Range("H122").Value = 605266
$ws.Range("I122").Value = 1136169.1
$ws.Range("K122").Value = 3408507.3
$ws.Range("M122").Value = -3406057.3
$ws.Range("H132").Value = 1958.3334
$ws.Range("I132").Value = 1958.3334
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5875.0002
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -3345.0002
$ws.Range("H134").Value = 3130.8276
$ws.Range("I134").Value = 2556.3333
$ws.Range("J134").Value = 4638.875
$ws.Range("K134").Value = 7668.999899999999
$ws.Range("L134").Value = 13916.625
$ws.Range("M134").Value = -5133.999899999999
$ws.Range("N134").Value = -18986.625
$ws.Range("H136").Value = 2349.32
$ws.Range("I136").Value = 1402.1666
$ws.Range("K136").Value = 4206.4998
$ws.Range("M136").Value = -1656.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 747.25
$ws.Range("I68").Value = 670.5
$ws.Range("K68").Value = 2011.5
$ws.Range("M68").Value = -1200.5
$ws.Range("H71").Value = 747.25
$ws.Range("I71").Value = 670.5
$ws.Range("K71").Value = 6034.5
$ws.Range("M71").Value = -1978.5
$ws.Range("H123").Value = 5000
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()
$ws.Range("H136").Value = 3099.4375
$ws.Range("J136").Value = 4090.7273
$ws.Range("L136").Value = 12272.1819
$ws.Range("N136").Value = -22472.1819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").ClearContents()
$ws.Range("N34").Value = 0
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").ClearContents()
$ws.Range("N76").Value = 0
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").ClearContents()
$ws.Range("N79").Value = 0
$ws.Range("H97").Value = 4443.7085
$ws.Range("I97").Value = 234.73685
$ws.Range("K97").Value = 234.73685
$ws.Range("M97").Value = 261.26315
$ws.Range("H102").Value = 10749
$ws.Range("I102").Value = 2500
$ws.Range("J102").Value = 18998
$ws.Range("K102").Value = 2500
$ws.Range("L102").Value = 18998
$ws.Range("M102").Value = -878
$ws.Range("N102").Value = -22242
$ws.Range("H122").Value = 4283.8335
$ws.Range("I122").Value = 2446.75
$ws.Range("J122").Value = 4808.7144
$ws.Range("K122").Value = 7340.25
$ws.Range("L122").Value = 14426.1432
$ws.Range("M122").Value = -4890.25
$ws.Range("N122").Value = -19326.1432
$ws.Range("H126").Value = 3175.1667
$ws.Range("I126").Value = 3218.3635
$ws.Range("K126").Value = 9655.0905
$ws.Range("M126").Value = -7185.0905
$ws.Range("H132").Value = 1517.8
$ws.Range("I132").Value = 1462.8485
$ws.Range("K132").Value = 4388.5455
$ws.Range("M132").Value = -1858.5455

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3757.6667
$ws.Range("I93").Value = 3757.6667
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 3757.6667
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -2509.6667
$ws.Range("H122").Value = 4857.533
$ws.Range("I122").Value = 2421.7273
$ws.Range("K122").Value = 7265.1819
$ws.Range("M122").Value = -4815.1819
$ws.Range("H136").Value = 1932.4445
$ws.Range("I136").Value = 1763.2051
$ws.Range("K136").Value = 5289.615299999999
$ws.Range("M136").Value = -2739.615299999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3610.1428
$ws.Range("I81").Value = 2608
$ws.Range("K81").Value = 5216
$ws.Range("M81").Value = -4155
$ws.Range("H84").Value = 3610.1428
$ws.Range("I84").Value = 2608
$ws.Range("K84").Value = 26080
$ws.Range("M84").Value = -20776
$ws.Range("H122").Value = 2800.6086
$ws.Range("I122").Value = 2810.45
$ws.Range("K122").Value = 8431.349999999999
$ws.Range("M122").Value = -5981.349999999999
